# Fixed Stimulus Absolute Timestamps
$wb = $excel.ActiveWorkbook

# Rename worksheets (order corresponds to sheetId 1..5 / rId1..rId5)
$wb.Worksheets.Item(1).Name = "GNG_TO-16504778717464905"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778737005248"
$wb.Worksheets.Item(3).Name = "RS_TO-1650477873706491"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504778737654963"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778738285253"

# Sheet 1 - GNG
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778717036676.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778717296362.csv"
$ws1.Range("B4").Value = "go_stims-1650477871731491.csv"
$ws1.Range("B5").Value = "GNG_stims-1650477871745524.csv"

# Sheet 2 - NB
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-165047787322149.csv"
$ws2.Range("B3").Value = "ZB-match_2-1650477872050497.csv"
$ws2.Range("B4").Value = "TB-16504778735515144.csv"
$ws2.Range("B5").Value = "OB-16504778729044964.csv"
$ws2.Range("B6").Value = "TB-16504778734575288.csv"
$ws2.Range("B7").Value = "TB-16504778736865246.csv"
$ws2.Range("B8").Value = "ZB-match_1-16504778722724934.csv"
$ws2.Range("B9").Value = "ZB-match_7-1650477871899526.csv"
$ws2.Range("B10").Value = "OB-16504778728025267.csv"

# Sheet 3 - RS
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# Sheet 4 - TOL
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778737325244.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778737084918.csv"
$ws4.Range("B4").Value = "MM_stims-16504778737485242.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778737334921.csv"
$ws4.Range("B6").Value = "MM_stims-16504778737644894.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778737494905.csv"

# Sheet 5 - vSAT
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16504778737965252.csv"
$ws5.Range("B3").Value = "SAT_stims-16504778737684903.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504778738125253.csv"
$ws5.Range("B5").Value = "SAT_stims-16504778737805257.csv"
